$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (D value or $null, E value)
# D values that are null mean the Price column text did not change in the
# source diff for that row - only Volume(1h) changed.
$updates = @{
    2  = @("59.891.03", "  -0.06%  ")
    3  = @("2.385.43", "  -1.29%  ")
    4  = @($null, "  -0.03%  ")
    5  = @("555.21", "  +0.61%  ")
    6  = @("133.37", "  -2.89%  ")
    7  = @($null, "  -0.01%  ")
    8  = @($null, "  -0.93%  ")
    10 = @("5.62", "  -1.53%  ")
    11 = @($null, "  +1.20%  ")
    12 = @($null, "  -3.04%  ")
    13 = @("24.43", "  -4.34%  ")
    14 = @("2.809.17", "  -1.33%  ")
    15 = @("59.795.32", "  -0.09%  ")
    16 = @("0.0000137", "  -0.79%  ")
    17 = @("2.386.50", "  -2.32%  ")
    18 = @($null, "  -1.88%  ")
    19 = @($null, "  +1.91%  ")
    20 = @("320.60", "  -2.61%  ")
    21 = @($null, "  +1.04%  ")
    22 = @($null, "  +0.02%  ")
    23 = @("64.16", "  -3.60%  ")
    24 = @("0.173", "  +0.37%  ")
    25 = @($null, "  +0.13%  ")
    26 = @("8.44", "  -2.44%  ")
    27 = @($null, "  +0.42%  ")
    28 = @($null, "  +1.66%  ")
    29 = @($null, "  -2.05%  ")
    30 = @("169.72", "  +0.86%  ")
    31 = @($null, "  -1.01%  ")
    32 = @("1.07", "  +5.73%  ")
    33 = @($null, "  -2.59%  ")
    34 = @("18.16", "  -2.49%  ")
    35 = @($null, "  +0.01%  ")
    36 = @($null, "  +1.25%  ")
    37 = @($null, "  +0.01%  ")
    38 = @($null, "  -2.33%  ")
    39 = @($null, "  -1.61%  ")
    40 = @("318.54", "  +1.35%  ")
    41 = @($null, "  -2.40%  ")
    42 = @("145.57", "  +4.89%  ")
    43 = @($null, "  -4.25%  ")
    44 = @("0.0968", "  +0.06%  ")
    45 = @("19.70", "  +0.81%  ")
    46 = @($null, "  -1.76%  ")
    47 = @("0.572", "  -1.49%  ")
    48 = @($null, "  -2.95%  ")
    49 = @($null, "  -0.07%  ")
    50 = @($null, "  +0.05%  ")
    51 = @("0.948", "  -0.08%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    if ($dVal -ne $null) {
        $dCell = $ws.Range("D$row")
        # These "Price" strings look numeric (e.g. "555.21"); Excel's COM
        # layer auto-converts plain assignment into a real number. Force
        # text storage via a Text number format, then strip the format
        # back off so the cell's style/appearance is left untouched -
        # only its stored value type (text) changes, matching the source
        # file where these are plain inline strings with no associated
        # style.
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.ClearFormats()
    }

    $ws.Range("E$row").Value = $eVal
}
